$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.037.07"
$ws.Range("D2").Style = $s
$ws.Range("E2").Value = "  -3.46%  "

$s = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.600.89"
$ws.Range("D3").Style = $s
$ws.Range("E3").Value = "  -2.79%  "

$s = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = $s
$ws.Range("E4").Value = "  +0.26%  "

$s = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.003"
$ws.Range("D5").Style = $s
$ws.Range("E5").Value = "  +0.08%  "

$s = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "300.72"
$ws.Range("D6").Style = $s
$ws.Range("E6").Value = "  -2.64%  "

$s = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3742"
$ws.Range("D7").Style = $s
$ws.Range("E7").Value = "  -3.88%  "

$s = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3605"
$ws.Range("D8").Style = $s
$ws.Range("E8").Value = "  -5.75%  "

$s = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.57"
$ws.Range("D9").Style = $s
$ws.Range("E9").Value = "  -5.48%  "

$s = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.006"
$ws.Range("D10").Style = $s
$ws.Range("E10").Value = "  +0.35%  "

$s = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.250"
$ws.Range("D11").Style = $s
$ws.Range("E11").Value = "  -7.04%  "

$s = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07974"
$ws.Range("D12").Style = $s
$ws.Range("E12").Value = "  -5.34%  "

$s = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.65"
$ws.Range("D13").Style = $s
$ws.Range("E13").Value = "  -4.98%  "

$s = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.509"
$ws.Range("D14").Style = $s
$ws.Range("E14").Value = "  -7.88%  "

$s = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.296"
$ws.Range("D15").Style = $s
$ws.Range("E15").Value = "  -7.48%  "

$s = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001231"
$ws.Range("D16").Style = $s
$ws.Range("E16").Value = "  -6.37%  "

$s = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.601.47"
$ws.Range("D17").Style = $s
$ws.Range("E17").Value = "  -2.90%  "

$s = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.93"
$ws.Range("D18").Style = $s
$ws.Range("E18").Value = "  -3.60%  "

$s = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06774"
$ws.Range("D19").Style = $s
$ws.Range("E19").Value = "  -2.91%  "

$s = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.15"
$ws.Range("D20").Style = $s
$ws.Range("E20").Value = "  -7.44%  "

$s = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.510"
$ws.Range("D21").Style = $s
$ws.Range("E21").Value = "  -5.97%  "

$s = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.005"
$ws.Range("D22").Style = $s
$ws.Range("E22").Value = "  +0.27%  "

$s = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.97"
$ws.Range("D23").Style = $s
$ws.Range("E23").Value = "  -5.04%  "

$s = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.058.31"
$ws.Range("D24").Style = $s
$ws.Range("E24").Value = "  -3.40%  "

$s = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.356"
$ws.Range("D25").Style = $s

$s = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.810"
$ws.Range("D26").Style = $s
$ws.Range("E26").Value = "  -5.04%  "

$s = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.94"
$ws.Range("D27").Style = $s
$ws.Range("E27").Value = "  -4.62%  "

$s = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.67"
$ws.Range("D28").Style = $s
$ws.Range("E28").Value = "  -0.74%  "

$s = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.258"
$ws.Range("D29").Style = $s
$ws.Range("E29").Value = "  -2.40%  "

$s = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.37"
$ws.Range("D30").Style = $s
$ws.Range("E30").Value = "  -4.25%  "

$s = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.398"
$ws.Range("D31").Style = $s
$ws.Range("E31").Value = "  -4.61%  "

$s = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.690"
$ws.Range("D32").Style = $s
$ws.Range("E32").Value = "  -13.88%  "

$s = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.778.80"
$ws.Range("D33").Style = $s
$ws.Range("E33").Value = "  -2.84%  "

$s = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9547"
$ws.Range("D34").Style = $s
$ws.Range("E34").Value = "  -8.53%  "

$s = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07636"
$ws.Range("D35").Style = $s
$ws.Range("E35").Value = "  -4.76%  "

$s = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02752"
$ws.Range("D36").Style = $s
$ws.Range("E36").Value = "  -6.67%  "

$s = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2523"
$ws.Range("D37").Style = $s
$ws.Range("E37").Value = "  -5.55%  "

$s = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.133"
$ws.Range("D38").Style = $s
$ws.Range("E38").Value = "  -8.41%  "

$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$s = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08812"
$ws.Range("D39").Style = $s
$ws.Range("E39").Value = "  -3.03%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$s = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.03"
$ws.Range("D40").Style = $s
$ws.Range("E40").Value = "  -7.52%  "

$s = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.386"
$ws.Range("D41").Style = $s
$ws.Range("E41").Value = "  -2.23%  "

$s = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7088"
$ws.Range("D42").Style = $s
$ws.Range("E42").Value = "  -6.10%  "

$s = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.64"
$ws.Range("D43").Style = $s
$ws.Range("E43").Value = "  -5.91%  "

$s = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.57"
$ws.Range("D44").Style = $s
$ws.Range("E44").Value = "  -4.89%  "

$s = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6552"
$ws.Range("D45").Style = $s
$ws.Range("E45").Value = "  -5.57%  "

$s = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.002"
$ws.Range("D46").Style = $s
$ws.Range("E46").Value = "  +0.06%  "

$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$s = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.966"
$ws.Range("D47").Style = $s
$ws.Range("E47").Value = "  -2.83%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$s = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.268"
$ws.Range("D48").Style = $s
$ws.Range("E48").Value = "  -7.51%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$s = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07976"
$ws.Range("D49").Style = $s
$ws.Range("E49").Value = "  -3.61%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$s = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.02"
$ws.Range("D50").Style = $s
$ws.Range("E50").Value = "  -2.18%  "

$s = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.161"
$ws.Range("D51").Style = $s
$ws.Range("E51").Value = "  -3.68%  "
